$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Telefone column (D) currently holds formatted phone-number strings such as
# "(69) 3819-7424". The fix stores the digits only, as real numbers, instead
# of formatted text (fixing the "field limit" bug from the commit message).
$ws.Range("D2").Value = 6938197424
$ws.Range("D3").Value = 8625418136
$ws.Range("D4").Value = 9236144316
$ws.Range("D5").Value = 8328964906

# Column widths were adjusted (re-autofit) after the edit.
$ws.Columns.Item(1).ColumnWidth = 23.83
$ws.Columns.Item(2).ColumnWidth = 28.83

# Restore the selected cell as recorded in the saved view state.
$ws.Range("E10").Select()
